$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# --- Add a new "Dijkstra:" results block in rows 66-72, mirroring the
# --- existing "MSDijkstra:" block (rows 53-59).

# Header row (row 66)
$ws.Range("A66").Value = "Dijkstra:"
$ws.Range("B66").Value = "PenafielScc:"
$ws.Range("C66").Value = "EspinhoScc:"
$ws.Range("D66").Value = "Penafiel:"
$ws.Range("E66").Value = "Espinho:"
$ws.Range("F66").Value = "PortoScc:"
$ws.Range("G66").Value = "Porto:"

# Data rows (67-71)
$ws.Range("B67").Value = 19.959
$ws.Range("C67").Value = 114.732
$ws.Range("D67").Value = 30.884699999999999
$ws.Range("E67").Value = 96.261200000000002
$ws.Range("F67").Value = 3693.31
$ws.Range("G67").Value = 2478.59

$ws.Range("B68").Value = 28.923999999999999
$ws.Range("C68").Value = 84.806899999999999
$ws.Range("D68").Value = 24.933199999999999
$ws.Range("E68").Value = 286.90800000000002
$ws.Range("F68").Value = 3638.23
$ws.Range("G68").Value = 3857.6

$ws.Range("B69").Value = 22.9391
$ws.Range("C69").Value = 103.42700000000001
$ws.Range("D69").Value = 19.947600000000001
$ws.Range("E69").Value = 101.904
$ws.Range("F69").Value = 3307.27
$ws.Range("G69").Value = 2505.56

$ws.Range("B70").Value = 24.932400000000001
$ws.Range("C70").Value = 101.52800000000001
$ws.Range("D70").Value = 24.965800000000002
$ws.Range("E70").Value = 124.667
$ws.Range("F70").Value = 4043.88
$ws.Range("G70").Value = 3041.3

$ws.Range("B71").Value = 29.9541
$ws.Range("C71").Value = 110.739
$ws.Range("D71").Value = 16.966000000000001
$ws.Range("E71").Value = 87.535300000000007
$ws.Range("F71").Value = 4716.62
$ws.Range("G71").Value = 4224.71

# Average row (72) - filled across B72:G72 in one go so the formula is
# shared across the row (mirrors how the other "Media:" rows were filled).
$ws.Range("A72").Value = "Media:"
$ws.Range("B72:G72").Formula = "=SUM(B67:B71)/5"

# --- Update the sheet view: selection moves to A67 (scroll position
# --- itself - topLeftCell - is not reproducible through this COM surface).
$ws.Range("A67").Select()
